$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.703100000000004
$ws.Range("E3").Value = 15.62760000000001
$ws.Range("B12").Value = 4.902099999999997
$ws.Range("D14").Value = -7.566299999999999
$ws.Range("E20").Value = 16.42109999999999
$ws.Range("E25").Value = 17.07920000000001
$ws.Range("D26").Value = -8.606600000000009
$ws.Range("B27").Value = 5.876400000000001
$ws.Range("E30").Value = 15.70640000000001
$ws.Range("D31").Value = -8.782799999999998
$ws.Range("B32").Value = 6.605100000000001
$ws.Range("D35").Value = -8.617800000000003
$ws.Range("B36").Value = 9.341500000000003
$ws.Range("D37").Value = -7.8175
$ws.Range("B38").Value = 4.8209
$ws.Range("E44").Value = 16.2783
$ws.Range("D45").Value = -7.5367
$ws.Range("B46").Value = 6.997800000000004
$ws.Range("E47").Value = 16.04819999999999
$ws.Range("D52").Value = -7.729199999999995
$ws.Range("B54").Value = 5.050299999999998
$ws.Range("B55").Value = 5.291499999999998
$ws.Range("B56").Value = 5.136099999999999
$ws.Range("D57").Value = -8.584899999999999
$ws.Range("E58").Value = 16.7531
$ws.Range("B67").Value = 5.789400000000003
$ws.Range("B69").Value = 5.178499999999998
$ws.Range("B72").Value = 4.991100000000004
$ws.Range("E78").Value = 16.64260000000003
$ws.Range("D81").Value = -7.123799999999996
$ws.Range("B83").Value = 5.389899999999998
$ws.Range("D83").Value = -9.162399999999995
$ws.Range("E84").Value = 16.3805
$ws.Range("B86").Value = 4.897600000000005
$ws.Range("E89").Value = 17.34380000000002
$ws.Range("B91").Value = 5.823800000000001
$ws.Range("E91").Value = 17.99540000000001
$ws.Range("E92").Value = 18.05320000000001
$ws.Range("B93").Value = 6.135100000000003
$ws.Range("E96").Value = 15.47099999999999
$ws.Range("B99").Value = 4.798399999999997
$ws.Range("D100").Value = -8.633900000000001
$ws.Range("D102").Value = -7.666799999999999
$ws.Range("E102").Value = 16.54859999999999
